$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column header "Household Size" in I1
$ws.Range("I1").Value = "Household Size"

# Household Size values for rows 2-11
$householdSizes = @(1, 2, 4, 3, 2, 1, 2, 2, 1, 1)
for ($i = 0; $i -lt $householdSizes.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $householdSizes[$i]
}

# Correct some Monthly Household Income values that also changed
$ws.Range("H4").Value = 2000
$ws.Range("H7").Value = 500
$ws.Range("H9").Value = 500
$ws.Range("H10").Value = 500

# Update the selection to match the post-edit state (I12)
$ws.Range("I12").Select()
